# Scheduled market-data refresh: updates Universalis price snapshots
# (currentAveragePrice/NQ/HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1165.56
$ws.Range("J18").Value = 1263
$ws.Range("L18").Value = 1263
$ws.Range("N18").Value = -1831
$ws.Range("H86").Value = 85011.086
$ws.Range("I86").Value = 112855.89
$ws.Range("J86").Value = 1476.6666
$ws.Range("K86").Value = 112855.89
$ws.Range("L86").Value = 1476.6666
$ws.Range("M86").Value = -111732.89
$ws.Range("N86").Value = -3722.6666
$ws.Range("H89").Value = 85011.086
$ws.Range("I89").Value = 112855.89
$ws.Range("J89").Value = 1476.6666
$ws.Range("K89").Value = 564279.45
$ws.Range("L89").Value = 7383.333000000001
$ws.Range("M89").Value = -558663.45
$ws.Range("N89").Value = -18615.333
$ws.Range("H92").Value = 930.8333
$ws.Range("I92").Value = 563.3333
$ws.Range("J92").Value = 2033.3334
$ws.Range("K92").Value = 563.3333
$ws.Range("L92").Value = 2033.3334
$ws.Range("M92").Value = 684.6667
$ws.Range("N92").Value = -4529.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 399.375
$ws.Range("I97").Value = 322.30768
$ws.Range("K97").Value = 322.30768
$ws.Range("M97").Value = 173.69232

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2351.2058
$ws.Range("I20").Value = 2469.8948
$ws.Range("K20").Value = 2469.8948
$ws.Range("M20").Value = -2222.8948
$ws.Range("H86").Value = 2821.9
$ws.Range("I86").Value = 2743.5715
$ws.Range("J86").Value = 3004.6667
$ws.Range("K86").Value = 2743.5715
$ws.Range("L86").Value = 3004.6667
$ws.Range("M86").Value = -1620.5715
$ws.Range("N86").Value = -5250.6667
$ws.Range("H89").Value = 2821.9
$ws.Range("I89").Value = 2743.5715
$ws.Range("J89").Value = 3004.6667
$ws.Range("K89").Value = 13717.8575
$ws.Range("L89").Value = 15023.3335
$ws.Range("M89").Value = -8101.8575
$ws.Range("N89").Value = -26255.3335
$ws.Range("H94").Value = 611.1539
$ws.Range("I94").Value = 565.625
$ws.Range("J94").Value = 684
$ws.Range("K94").Value = 565.625
$ws.Range("L94").Value = 684
$ws.Range("M94").Value = -114.625
$ws.Range("N94").Value = -1586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 500
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 500
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -726
$ws.Range("H41").Value = 8033.3335
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 8033.3335
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 8033.3335
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -8889.333500000001
$ws.Range("H50").Value = 10280
$ws.Range("J50").Value = 10280
$ws.Range("L50").Value = 10280
$ws.Range("N50").Value = -11530
$ws.Range("H51").Value = 13781.272
$ws.Range("J51").Value = 13781.272
$ws.Range("L51").Value = 13781.272
$ws.Range("N51").Value = -15253.272
$ws.Range("H59").Value = 44400
$ws.Range("J59").Value = 44400
$ws.Range("L59").Value = 44400
$ws.Range("N59").Value = -46690
$ws.Range("H60").Value = 14087.429
$ws.Range("J60").Value = 14501.846
$ws.Range("L60").Value = 14501.846
$ws.Range("N60").Value = -15523.846
$ws.Range("H61").Value = 13781.272
$ws.Range("J61").Value = 13781.272
$ws.Range("L61").Value = 13781.272
$ws.Range("N61").Value = -14477.272
$ws.Range("H62").Value = 2305
$ws.Range("J62").Value = 2200
$ws.Range("L62").Value = 2200
$ws.Range("N62").Value = -3448
$ws.Range("H65").Value = 2305
$ws.Range("J65").Value = 2200
$ws.Range("L65").Value = 11000
$ws.Range("N65").Value = -17240
$ws.Range("H107").Value = 540.931
$ws.Range("I107").Value = 534.7895
$ws.Range("J107").Value = 552.6
$ws.Range("K107").Value = 534.7895
$ws.Range("L107").Value = 552.6
$ws.Range("M107").Value = 1385.2105
$ws.Range("N107").Value = -4392.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 80
$ws.Range("I15").Value = 80
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 240
$ws.Range("L15").Value = 240
$ws.Range("M15").Value = -100
$ws.Range("N15").Value = -520
$ws.Range("H62").Value = 3277.8462
$ws.Range("J62").Value = 3392.6667
$ws.Range("L62").Value = 10178.0001
$ws.Range("N62").Value = -11550.0001
$ws.Range("H65").Value = 3277.8462
$ws.Range("J65").Value = 3392.6667
$ws.Range("L65").Value = 30534.0003
$ws.Range("N65").Value = -37398.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7777.6665
$ws.Range("I80").Value = 6000
$ws.Range("J80").Value = 7999.875
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 7999.875
$ws.Range("M80").Value = -5002
$ws.Range("N80").Value = -9995.875
$ws.Range("H83").Value = 7777.6665
$ws.Range("I83").Value = 6000
$ws.Range("J83").Value = 7999.875
$ws.Range("K83").Value = 30000
$ws.Range("L83").Value = 39999.375
$ws.Range("M83").Value = -25008
$ws.Range("N83").Value = -49983.375
$ws.Range("H97").Value = 901.8570999999999
$ws.Range("I97").Value = 922.1429000000001
$ws.Range("J97").Value = 861.2857
$ws.Range("K97").Value = 922.1429000000001
$ws.Range("L97").Value = 861.2857
$ws.Range("M97").Value = -426.1429000000001
$ws.Range("N97").Value = -1853.2857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1742.875
$ws.Range("I68").Value = 2008
$ws.Range("J68").Value = 1301
$ws.Range("K68").Value = 2008
$ws.Range("L68").Value = 1301
$ws.Range("M68").Value = -1259
$ws.Range("N68").Value = -2799
$ws.Range("H71").Value = 1742.875
$ws.Range("I71").Value = 2008
$ws.Range("J71").Value = 1301
$ws.Range("K71").Value = 10040
$ws.Range("L71").Value = 1301
$ws.Range("M71").Value = -6296
$ws.Range("N71").Value = -13993
$ws.Range("H93").Value = 3380577.8
$ws.Range("I93").Value = 3380577.8
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 3380577.8
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -3379329.8
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 775
$ws.Range("I100").Value = 962.5
$ws.Range("J100").Value = 400
$ws.Range("K100").Value = 962.5
$ws.Range("L100").Value = 400
$ws.Range("M100").Value = -421.5
$ws.Range("N100").Value = -1482

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 674.4706
$ws.Range("I113").Value = 748.8570999999999
$ws.Range("J113").Value = 327.33334
$ws.Range("K113").Value = 2246.5713
$ws.Range("L113").Value = 982.0000200000001
$ws.Range("M113").Value = -76.57129999999961
$ws.Range("N113").Value = -5322.00002
$ws.Range("H126").Value = 13322.111
$ws.Range("I126").Value = 16557.715
$ws.Range("J126").Value = 1997.5
$ws.Range("K126").Value = 49673.145
$ws.Range("L126").Value = 5992.5
$ws.Range("M126").Value = -47203.145
$ws.Range("N126").Value = -10932.5
